# Update "想去人数" (F column) counts on the 展览 sheet (index 1)
# and the matching rows on the 全部类型 sheet (index 4), reflecting the
# refreshed scrape data for "output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# Row -> new F-column value, for the 展览 sheet
$exhibitUpdates = @{
    5  = 271
    6  = 1075
    7  = 1406
    10 = 737
    13 = 122
    14 = 415
    15 = 1308
    17 = 86
    18 = 265
    20 = 639
    24 = 5618
    26 = 116
    29 = 14148
    30 = 1410
    32 = 89
    35 = 578
    36 = 4166
    37 = 107
    38 = 352
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Row -> new F-column value, for the 全部类型 sheet (same events, shifted
# row numbers because that sheet interleaves other categories)
$allUpdates = @{
    5  = 271
    6  = 1075
    7  = 1406
    10 = 737
    13 = 122
    14 = 415
    15 = 1308
    17 = 86
    18 = 265
    21 = 639
    27 = 5618
    29 = 116
    32 = 14148
    33 = 1410
    35 = 89
    38 = 578
    39 = 4166
    40 = 107
    41 = 352
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
